# chilli-pop (Version 2): move the "Meta description" blurb from the top
# of the document down to the bottom (as a new bold heading-style line,
# without the literal "Meta description" label), and replace the old
# AI image-prompt paragraph's text with the blurb text, keeping that
# paragraph's italic run formatting intact.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the whole "Meta description: ..." paragraph that currently
#    sits right after the H1 title at the top of the document.
# ---------------------------------------------------------------------
$d.Paragraphs.Item(2).Range.Delete()

# ---------------------------------------------------------------------
# 2) Insert a brand-new paragraph right before the very last paragraph
#    (the old AI image-generation prompt) containing a leading empty
#    run followed by a bold run with the title text. We splice in a
#    tiny two-paragraph OOXML fragment so Word creates a real
#    paragraph break with clean (inherited-free) formatting, then glue
#    the trailing placeholder paragraph back onto the original last
#    paragraph's content.
# ---------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$xmlSnippet = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Chilli Pop for Free – A Mexican Cuisine-Inspired Cluster-Based Game</w:t></w:r></w:p><w:p><w:r><w:t>ZZZ_TMP_SPLIT_ZZZ</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($xmlSnippet)

# Drop the temporary placeholder text introduced purely to force the
# paragraph split above; the original final paragraph's own runs
# (leading empty run + italic run) are left completely untouched.
$d.Content.Find.Execute("ZZZ_TMP_SPLIT_ZZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ---------------------------------------------------------------------
# 3) Swap the old AI image prompt text for the meta-description text in
#    what is now the last paragraph, preserving its italic run.
# ---------------------------------------------------------------------
$oldText = "Create a vibrant feature image for Chilli Pop that showcases the game's Mexican theme and fun characters. The image should be in a cartoon style and feature a happy Maya warrior with glasses as the main focus. Surround the warrior with symbols from the game, such as tomatoes, garlic, peppers, and onions that have been transformed into wacky characters. Use bright colors to make the image pop and ensure that it captures the spirit of the game's exciting gameplay and cluster-based winning combinations. Add the Chilli Pop logo to the center of the image to tie it all together and make it clear which game it represents."
$newText = "Read our review of Chilli Pop to learn more about its gameplay mechanics, payouts, and features. Play for free and enjoy the Mexican cuisine-inspired theme."
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
